$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.011.52"
$ws.Range("E2").Value = "'  +3.40%  "
$ws.Range("D3").Value = "'2.418.96"
$ws.Range("E3").Value = "'  +3.02%  "
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("D5").Value = "'552.73"
$ws.Range("E5").Value = "'  +2.04%  "
$ws.Range("D6").Value = "'137.16"
$ws.Range("E6").Value = "'  +2.22%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'0.583"
$ws.Range("E8").Value = "'  +2.96%  "
$ws.Range("D9").Value = "'0.107"
$ws.Range("E9").Value = "'  +2.53%  "
$ws.Range("E10").Value = "'  +5.33%  "
$ws.Range("D11").Value = "'0.359"
$ws.Range("E11").Value = "'  +0.65%  "
$ws.Range("E12").Value = "'  -1.92%  "
$ws.Range("D13").Value = "'24.88"
$ws.Range("E13").Value = "'  +4.69%  "
$ws.Range("D14").Value = "'2.851.71"
$ws.Range("D15").Value = "'59.970.38"
$ws.Range("E15").Value = "'  +3.46%  "
$ws.Range("D17").Value = "'2.422.87"
$ws.Range("E17").Value = "'  +3.65%  "
$ws.Range("D18").Value = "'11.41"
$ws.Range("E18").Value = "'  +6.29%  "
$ws.Range("E19").Value = "'  +2.15%  "
$ws.Range("D20").Value = "'331.82"
$ws.Range("E20").Value = "'  +0.95%  "
$ws.Range("E21").Value = "'  -0.06%  "
$ws.Range("E22").Value = "'  -0.12%  "
$ws.Range("D23").Value = "'65.48"
$ws.Range("E23").Value = "'  +3.77%  "
$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "'  +3.61%  "
$ws.Range("D25").Value = "'8.62"
$ws.Range("E25").Value = "'  +2.97%  "
$ws.Range("E26").Value = "'  +0.54%  "
$ws.Range("E27").Value = "'  +1.20%  "
$ws.Range("D28").Value = "'0.0₃0782"
$ws.Range("E28").Value = "'  +6.15%  "
$ws.Range("E29").Value = "'  +0.80%  "
$ws.Range("D30").Value = "'170.61"
$ws.Range("E30").Value = "'  +0.15%  "
$ws.Range("D31").Value = "'6.23"
$ws.Range("E31").Value = "'  +1.42%  "
$ws.Range("D32").Value = "'18.61"
$ws.Range("E32").Value = "'  +1.55%  "
$ws.Range("E33").Value = "'  +1.87%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "'  -0.03%  "
$ws.Range("E35").Value = "'  +5.36%  "
$ws.Range("D37").Value = "'4.20"
$ws.Range("E37").Value = "'  +0.12%  "
$ws.Range("E38").Value = "'  +0.64%  "
$ws.Range("D39").Value = "'39.59"
$ws.Range("E39").Value = "'  +1.18%  "
$ws.Range("E40").Value = "'  +9.79%  "
$ws.Range("D41").Value = "'313.70"
$ws.Range("E41").Value = "'  +8.14%  "
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "'  +1.39%  "
$ws.Range("D43").Value = "'139.53"
$ws.Range("E43").Value = "'  -0.72%  "
$ws.Range("D44").Value = "'0.0964"
$ws.Range("E44").Value = "'  +1.32%  "
$ws.Range("E45").Value = "'  +1.80%  "
$ws.Range("D46").Value = "'19.40"
$ws.Range("E46").Value = "'  +2.61%  "
$ws.Range("B47").Value = "'Mantle"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.575"
$ws.Range("E47").Value = "'  +1.75%  "
$ws.Range("B48").Value = "'Polygon"
$ws.Range("C48").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.410"
$ws.Range("E48").Value = "'  +7.37%  "
$ws.Range("E49").Value = "'  +1.65%  "
$ws.Range("E50").Value = "'  +1.35%  "
$ws.Range("D51").Value = "'11.05"
$ws.Range("E51").Value = "'  -0.16%  "
